$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the existing header style (copy format from F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Update existing MSE/R2/MAE values
$ws.Range("B2").Value = 0.7331181174200663
$ws.Range("C2").Value = 0.9784664004414199
$ws.Range("D2").Value = 0.6623029683249402

$ws.Range("B3").Value = 0.2410320860048795
$ws.Range("C3").Value = 0.9966688713085116
$ws.Range("D3").Value = 0.3680431895364141

$ws.Range("B4").Value = 0.1542566699619067
$ws.Range("C4").Value = 0.9984101093943154
$ws.Range("D4").Value = 0.3357392818625894

$ws.Range("B5").Value = 0.4045655439972058
$ws.Range("C5").Value = 0.9975901184861534
$ws.Range("D5").Value = 0.5182038504872997

# Add new Elapsed Time / CPU values for each data row
$ws.Range("G2").Value = 0.4794827245333484
$ws.Range("H2").Value = 0.996

$ws.Range("G3").Value = 0.4794827245333484
$ws.Range("H3").Value = 0.996

$ws.Range("G4").Value = 0.4794827245333484
$ws.Range("H4").Value = 0.996

$ws.Range("G5").Value = 0.4794827245333484
$ws.Range("H5").Value = 0.996

Write-Output "done"
